$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the first "Completed Dt" value with the date number format (built-in
# format 14, "mm-dd-yy"), then propagate that same format to the rest of the
# column by copying it, so every cell shares a single style entry.
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Value = (Get-Date -Year 2024 -Month 8 -Day 31).Date

$ws.Range("F2").Copy()
$ws.Range("F3:F16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the rest of the "Completed Dt" values for rows 3-16
$ws.Range("F3").Value  = (Get-Date -Year 2024 -Month 8 -Day 31).Date
$ws.Range("F4").Value  = (Get-Date -Year 2024 -Month 8 -Day 31).Date
$ws.Range("F5").Value  = (Get-Date -Year 2024 -Month 8 -Day 30).Date
$ws.Range("F6").Value  = (Get-Date -Year 2024 -Month 8 -Day 30).Date
$ws.Range("F7").Value  = (Get-Date -Year 2024 -Month 8 -Day 30).Date
$ws.Range("F8").Value  = (Get-Date -Year 2024 -Month 8 -Day 29).Date
$ws.Range("F9").Value  = (Get-Date -Year 2024 -Month 8 -Day 29).Date
$ws.Range("F10").Value = (Get-Date -Year 2024 -Month 8 -Day 29).Date
$ws.Range("F11").Value = (Get-Date -Year 2024 -Month 8 -Day 29).Date
$ws.Range("F12").Value = (Get-Date -Year 2024 -Month 8 -Day 28).Date
$ws.Range("F13").Value = (Get-Date -Year 2024 -Month 8 -Day 28).Date
$ws.Range("F14").Value = (Get-Date -Year 2024 -Month 8 -Day 28).Date
$ws.Range("F15").Value = (Get-Date -Year 2024 -Month 8 -Day 28).Date
$ws.Range("F16").Value = (Get-Date -Year 2024 -Month 8 -Day 28).Date

# Leave the selection where the user ended up after entering the data
$ws.Range("F12:F16").Select()
